# Weekly update: a new "Cebollín" price record for Terminal Hortofrutícola
# Agro Chillán was published. Insert it as a new row 26 (pushing the
# existing rows 26-57 down to 27-58) and populate it with the new values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new row at position 26; everything below shifts down one row.
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row with the new weekly record.
$ws.Cells.Item(26, 1).Value = 7
$ws.Cells.Item(26, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(26, 3).Value = "Ñuble"
$ws.Cells.Item(26, 4).Value = 44902
$ws.Cells.Item(26, 5).Value = 16
$ws.Cells.Item(26, 6).Value = 100112037
$ws.Cells.Item(26, 7).Value = "Cebollín"
$ws.Cells.Item(26, 8).Value = "Sin especificar"
$ws.Cells.Item(26, 9).Value = "Segunda"
$ws.Cells.Item(26, 10).Value = 200
$ws.Cells.Item(26, 11).Value = 500
$ws.Cells.Item(26, 12).Value = 500
$ws.Cells.Item(26, 13).Value = 500
$ws.Cells.Item(26, 14).Value = "$/paquete 6 unidades"
$ws.Cells.Item(26, 15).Value = "Provincia de Diguillín"
$ws.Cells.Item(26, 16).Value = 83
$ws.Cells.Item(26, 17).Value = 6
$ws.Cells.Item(26, 18).Value = "Hortaliza"
